$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (C, D, F, H change; A, B, E, G stay the same) ---
$ws.Columns.Item(3).ColumnWidth = 135.17   # C: 50 -> 136
$ws.Columns.Item(4).ColumnWidth = 37.17    # D: 51 -> 38
$ws.Columns.Item(6).ColumnWidth = 16.17    # F: 16 -> 17
$ws.Columns.Item(8).ColumnWidth = 36.17    # H: 29 -> 37

# Row 2
$ws.Range("A2").Value = "'1326897"
$ws.Range("B2").Value = 'https://aiesec.org/opportunity/global-talent/1326897'
$ws.Range("C2").Value = 'Global Marketing Immersion Program Mastering the Fundamentals of Web Advertising'
$ws.Range("D2").Value = '日本、愛知県名古屋市'
$ws.Range("E2").Value = 'No'
$ws.Range("F2").Value = '0 applicants'
$ws.Range("G2").Value = '9 - 12 Weeks'
$ws.Range("H2").Value = 'QUARTET COMMUNICATIONS Co.,Ltd.'

# Row 3
$ws.Range("A3").Value = "'1326896"
$ws.Range("B3").Value = 'https://aiesec.org/opportunity/global-talent/1326896'
$ws.Range("C3").Value = 'Crisis Designer'
$ws.Range("D3").Value = 'Londres, Reino Unido'
$ws.Range("E3").Value = 'Yes'
$ws.Range("F3").Value = '1 applicant'
$ws.Range("G3").Value = '6 - 18 Months'
$ws.Range("H3").Value = 'International SOS'

# Row 4
$ws.Range("A4").Value = "'1326893"
$ws.Range("B4").Value = 'https://aiesec.org/opportunity/global-talent/1326893'
$ws.Range("C4").Value = 'Financial Performance and Strategy Analyst'
$ws.Range("D4").Value = 'Mexico City, CDMX, Mexico'
$ws.Range("E4").Value = 'No'
$ws.Range("F4").Value = '2 applicants'
$ws.Range("G4").Value = '6 - 18 Months'
$ws.Range("H4").Value = 'Sodexo Mexico'

# Row 5
$ws.Range("A5").Value = "'1326892"
$ws.Range("B5").Value = 'https://aiesec.org/opportunity/global-talent/1326892'
$ws.Range("C5").Value = 'Business and Market Strategy Analyst'
$ws.Range("D5").Value = 'Mexico City, CDMX, Mexico'
$ws.Range("E5").Value = 'No'
$ws.Range("F5").Value = '2 applicants'
$ws.Range("G5").Value = '6 - 18 Months'
$ws.Range("H5").Value = 'Sodexo Mexico'

# Row 6
$ws.Range("A6").Value = "'1326891"
$ws.Range("B6").Value = 'https://aiesec.org/opportunity/global-talent/1326891'
$ws.Range("C6").Value = 'Experience in learning about the production management system of the plating line and creating  proposals for efficiency improvement.'
$ws.Range("D6").Value = '日本、愛知県名古屋市'
$ws.Range("E6").Value = 'No'
$ws.Range("F6").Value = '0 applicants'
$ws.Range("G6").Value = '9 - 12 Weeks'
$ws.Range("H6").Value = 'YADAGAWA Electric Plating Co.,Ltd.'

# Row 7
$ws.Range("A7").Value = "'1326890"
$ws.Range("B7").Value = 'https://aiesec.org/opportunity/global-talent/1326890'
$ws.Range("C7").Value = 'Occupational Health and Safety Projects Specialist'
$ws.Range("D7").Value = 'Mexico City, CDMX, Mexico'
$ws.Range("E7").Value = 'No'
$ws.Range("F7").Value = '0 applicants'
$ws.Range("G7").Value = '6 - 18 Months'
$ws.Range("H7").Value = 'Sodexo Mexico'

# Row 8
$ws.Range("A8").Value = "'1326889"
$ws.Range("B8").Value = 'https://aiesec.org/opportunity/global-talent/1326889'
$ws.Range("C8").Value = 'Experience in developing software systems and practicing marketing for the social implementation of drones'
$ws.Range("D8").Value = '日本、東京都東京'
$ws.Range("E8").Value = 'No'
$ws.Range("F8").Value = '0 applicants'
$ws.Range("G8").Value = '9 - 12 Weeks'
$ws.Range("H8").Value = 'ALL NIPPON AIRWAYS CO., LTD.'

# Row 9
$ws.Range("A9").Value = "'1326872"
$ws.Range("B9").Value = 'https://aiesec.org/opportunity/global-talent/1326872'
$ws.Range("C9").Value = 'Conduct market research and develop a deployment strategy to bring the online IBDP to more countries and regions'
$ws.Range("D9").Value = '日本、東京都東京'
$ws.Range("E9").Value = 'No'
$ws.Range("F9").Value = '1 applicant'
$ws.Range("G9").Value = '9 - 12 Weeks'
$ws.Range("H9").Value = 'Aoba-BBT, Inc.'

# Row 10
$ws.Range("A10").Value = "'1326869"
$ws.Range("B10").Value = 'https://aiesec.org/opportunity/global-talent/1326869'
$ws.Range("C10").Value = 'Demonstrative Virtual Design of AI English Learning Tools'
$ws.Range("D10").Value = '日本、東京都東京'
$ws.Range("E10").Value = 'No'
$ws.Range("F10").Value = '2 applicants'
$ws.Range("G10").Value = '9 - 12 Weeks'
$ws.Range("H10").Value = 'Aoba-BBT, Inc.'

# Row 11
$ws.Range("A11").Value = "'1326868"
$ws.Range("B11").Value = 'https://aiesec.org/opportunity/global-talent/1326868'
$ws.Range("C11").Value = 'Practical Experience in Tourism-Based Marketing and Market Research for International Service Promotion'
$ws.Range("D11").Value = '日本、東京都台東区'
$ws.Range("E11").Value = 'No'
$ws.Range("F11").Value = '0 applicants'
$ws.Range("G11").Value = '9 - 12 Weeks'
$ws.Range("H11").Value = 'Japan Dream Tour Co., Ltd.'

# Row 12
$ws.Range("A12").Value = "'1326867"
$ws.Range("B12").Value = 'https://aiesec.org/opportunity/global-talent/1326867'
$ws.Range("C12").Value = 'Hands-on experience in back-end and front-end development of image recognition AI'
$ws.Range("D12").Value = '日本、東京都東京'
$ws.Range("E12").Value = 'No'
$ws.Range("F12").Value = '0 applicants'
$ws.Range("G12").Value = '9 - 12 Weeks'
$ws.Range("H12").Value = 'Systems Nakashima Co., Ltd.'

# Row 13
$ws.Range("A13").Value = "'1326860"
$ws.Range("B13").Value = 'https://aiesec.org/opportunity/global-talent/1326860'
$ws.Range("C13").Value = 'Exploring the Healthcare Markets of Japan and India and Developing New Business Ideas within Non-Practical Training'
$ws.Range("D13").Value = '日本、東京都東京'
$ws.Range("E13").Value = 'No'
$ws.Range("F13").Value = '1 applicant'
$ws.Range("G13").Value = '9 - 12 Weeks'
$ws.Range("H13").Value = 'SECOM Medical System Co., Ltd.'

# Row 14
$ws.Range("A14").Value = "'1326843"
$ws.Range("B14").Value = 'https://aiesec.org/opportunity/global-talent/1326843'
$ws.Range("C14").Value = 'Software Developer'
$ws.Range("D14").Value = 'Delft, Nederland'
$ws.Range("E14").Value = 'No'
$ws.Range("F14").Value = '8 applicants'
$ws.Range("G14").Value = '6 - 18 Months'
$ws.Range("H14").Value = 'IRM systems'

# Row 15
$ws.Range("A15").Value = "'1326574"
$ws.Range("B15").Value = 'https://aiesec.org/opportunity/global-talent/1326574'
$ws.Range("C15").Value = 'Global Portfolio Specialist (Product Line Care) - Global Electrolux Talent Program'
$ws.Range("D15").Value = 'Stokholm, İsveç'
$ws.Range("E15").Value = 'Yes'
$ws.Range("F15").Value = '138 applicants'
$ws.Range("G15").Value = '6 - 18 Months'
$ws.Range("H15").Value = 'AB Electrolux'

# Row 16
$ws.Range("A16").Value = "'1323468"
$ws.Range("B16").Value = 'https://aiesec.org/opportunity/global-talent/1323468'
$ws.Range("C16").Value = 'Sales Account Manager'
$ws.Range("D16").Value = 'Cyberjaya, Selangor, Malaysia'
$ws.Range("E16").Value = 'No'
$ws.Range("F16").Value = '36 applicants'
$ws.Range("G16").Value = '6 - 18 Months'
$ws.Range("H16").Value = 'IX Telecom Sdn Bhd'

# Row 17
$ws.Range("A17").Value = "'1320223"
$ws.Range("B17").Value = 'https://aiesec.org/opportunity/global-talent/1320223'
$ws.Range("C17").Value = 'Corporate Communications Associate'
$ws.Range("D17").Value = 'Panamá, Provincia de Panamá, Panamá'
$ws.Range("E17").Value = 'No'
$ws.Range("F17").Value = '118 applicants'
$ws.Range("G17").Value = '6 - 18 Months'
$ws.Range("H17").Value = 'Grünenthal, S.A.'

# Row 18
$ws.Range("A18").Value = "'1316723"
$ws.Range("B18").Value = 'https://aiesec.org/opportunity/global-talent/1316723'
$ws.Range("C18").Value = '[Impact Porto Alegre]- Social Media'
$ws.Range("D18").Value = 'Porto Alegre, RS, Brasil'
$ws.Range("E18").Value = 'No'
$ws.Range("F18").Value = '79 applicants'
$ws.Range("G18").Value = '9 - 12 Weeks'
$ws.Range("H18").Value = 'ESCOLA GIORDANO BRUNO LTDA'

# Column A held numeric-looking text; strip the quote-prefix marker
# Excel leaves behind so the cells come back out as plain shared-string
# text (no leftover cell style) just like the original file.
$ws.Range("A2:A18").ClearFormats()

# --- Apply the "Premium = Yes" highlight style to the new Yes row (row 15) ---
# Row 3 already carries this style from the template; copy its formatting
# (fill/border/font) onto E15 without disturbing the shared style table.
$ws.Range("E3").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
